$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.702.41'
$ws.Cells.Item(2, 5).Value = '  -0.42%  '
$ws.Cells.Item(3, 4).Value = '3.911.31'
$ws.Cells.Item(3, 5).Value = '  +4.42%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '603.04'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.15%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '165.12'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.26%  '
$ws.Cells.Item(7, 4).Value = '3.908.55'
$ws.Cells.Item(7, 5).Value = '  +4.40%  '
$ws.Cells.Item(8, 5).Value = '  -0.28%  '
$ws.Cells.Item(9, 5).Value = '  -1.56%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.165'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.47%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.37'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.08%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.461'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.62%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '37.02'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.57%  '
$ws.Cells.Item(14, 5).Value = '  -0.80%  '
$ws.Cells.Item(15, 4).Value = '4.567.96'
$ws.Cells.Item(15, 5).Value = '  +4.47%  '
$ws.Cells.Item(16, 4).Value = '3.909.29'
$ws.Cells.Item(16, 5).Value = '  +4.02%  '
$ws.Cells.Item(17, 4).Value = '68.880.62'
$ws.Cells.Item(17, 5).Value = '  -0.07%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.44'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.27%  '
$ws.Cells.Item(19, 5).Value = '  -0.68%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.99'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -4.64%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.10'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -1.66%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '485.18'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.90%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.719'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.63%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000168'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +12.87%  '
$ws.Cells.Item(25, 5).Value = '  -0.15%  '
$ws.Cells.Item(26, 5).Value = '  -0.75%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.05'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.75%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.09'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.66%  '
$ws.Cells.Item(29, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.94'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.78%  '
$ws.Cells.Item(31, 4).Value = '4.063.42'
$ws.Cells.Item(31, 5).Value = '  +4.43%  '
$ws.Cells.Item(32, 2).Value = 'NEARProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.84'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.53%  '
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.38'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.80%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '31.98'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.40%  '
$ws.Cells.Item(35, 4).Value = '3.853.93'
$ws.Cells.Item(35, 5).Value = '  +4.59%  '
$ws.Cells.Item(36, 5).Value = '  -0.47%  '
$ws.Cells.Item(38, 5).Value = '  +1.34%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.88'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.68%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.09%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.319'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -1.49%  '
$ws.Cells.Item(42, 2).Value = 'Bittensor'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '437.07'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +3.20%  '
$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.99'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -2.83%  '
$ws.Cells.Item(44, 2).Value = 'OKB'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '48.45'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.16%  '
$ws.Cells.Item(45, 2).Value = 'Stacks'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.99'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -0.02%  '
$ws.Cells.Item(46, 2).Value = 'USDe'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.02%  '
$ws.Cells.Item(47, 2).Value = 'Cosmos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.47'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.66%  '
$ws.Cells.Item(48, 4).Value = '2.837.85'
$ws.Cells.Item(48, 5).Value = '  +2.08%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.27'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +10.13%  '
$ws.Cells.Item(50, 2).Value = 'Monero'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '142.21'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.44%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0355'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.75%  '
